$wb = $excel.ActiveWorkbook

# Remove the two extra/unused sheets ("Naveen1" and "Naveen2"); keep "login" and "Naveen".
$wb.Worksheets("Naveen1").Delete() | Out-Null
$wb.Worksheets("Naveen2").Delete() | Out-Null

# Work on the "login" sheet (the one with the actual data / formula bug).
$ws = $wb.Worksheets("login")
$ws.Activate() | Out-Null

# A2 was "ppp" -> now "wsss".
$ws.Range("A2").Value = "wsss"

# A3 was "kkk" -> now the text value "99999" (kept as text via a leading quote,
# i.e. entered with the "quote prefix" so Excel treats it as text, not a number).
$ws.Range("A3").Value = "'99999"

# B3 was the literal string "tom123" -> now a real formula that sums to 50
# (this is the actual "FORMULA cell value" bug fix referenced by the commit).
$ws.Range("B3").Formula = "=SUM(10+40)"

# Carry the quote-prefix text formatting from A3 over to B3 as well, matching
# the sibling cell's style (both A3 and B3 end up on the same quote-prefixed style).
$ws.Range("A3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

# Final selection ends up on B4.
$ws.Range("B4").Select() | Out-Null
